$wb = $excel.ActiveWorkbook

# --- ALC sheet (sheet index 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3900
$ws.Cells.Item(76, 9).Value = 3200
$ws.Cells.Item(76, 11).Value = 3200
$ws.Cells.Item(76, 13).Value = -2885
$ws.Cells.Item(79, 8).Value = 3900
$ws.Cells.Item(79, 9).Value = 3200
$ws.Cells.Item(79, 11).Value = 3200
$ws.Cells.Item(79, 13).Value = -2108
$ws.Cells.Item(86, 8).Value = 2614.3809
$ws.Cells.Item(86, 9).Value = 2550.2
$ws.Cells.Item(86, 10).Value = 2672.7273
$ws.Cells.Item(86, 11).Value = 2550.2
$ws.Cells.Item(86, 12).Value = 2672.7273
$ws.Cells.Item(86, 13).Value = -1427.2
$ws.Cells.Item(86, 14).Value = -4918.7273
$ws.Cells.Item(89, 8).Value = 2614.3809
$ws.Cells.Item(89, 9).Value = 2550.2
$ws.Cells.Item(89, 10).Value = 2672.7273
$ws.Cells.Item(89, 11).Value = 12751
$ws.Cells.Item(89, 12).Value = 13363.6365
$ws.Cells.Item(89, 13).Value = -7135
$ws.Cells.Item(89, 14).Value = -24595.6365
$ws.Cells.Item(129, 8).Value = 935.14545
$ws.Cells.Item(129, 10).Value = 1061.5581
$ws.Cells.Item(129, 12).Value = 3184.6743
$ws.Cells.Item(129, 14).Value = -13184.6743
$ws.Cells.Item(137, 8).Value = 1727.9796
$ws.Cells.Item(137, 9).Value = 1465.7222
$ws.Cells.Item(137, 10).Value = 2454.2307
$ws.Cells.Item(137, 11).Value = 4397.1666
$ws.Cells.Item(137, 12).Value = 7362.6921
$ws.Cells.Item(137, 13).Value = -1847.1666
$ws.Cells.Item(137, 14).Value = -12462.6921
$ws.Cells.Item(138, 8).Value = 1656.6061
$ws.Cells.Item(138, 9).Value = 952.1539
$ws.Cells.Item(138, 10).Value = 2436
$ws.Cells.Item(138, 11).Value = 2856.4617
$ws.Cells.Item(138, 12).Value = 7308
$ws.Cells.Item(138, 13).Value = 2283.5383
$ws.Cells.Item(138, 14).Value = -17588

# --- ARM sheet (sheet index 2) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11068.29
$ws.Cells.Item(32, 9).Value = 9437.322
$ws.Cells.Item(32, 10).Value = 25747
$ws.Cells.Item(32, 11).Value = 9437.322
$ws.Cells.Item(32, 12).Value = 25747
$ws.Cells.Item(32, 13).Value = -9150.322
$ws.Cells.Item(32, 14).Value = -26321
$ws.Cells.Item(45, 8).Value = 1143.75
$ws.Cells.Item(45, 9).Value = 1187.5
$ws.Cells.Item(45, 10).Value = 1100
$ws.Cells.Item(45, 11).Value = 1187.5
$ws.Cells.Item(45, 12).Value = 1100
$ws.Cells.Item(45, 13).Value = -810.5
$ws.Cells.Item(45, 14).Value = -1854
$ws.Cells.Item(121, 8).Value = 41995
$ws.Cells.Item(121, 10).Value = 41995
$ws.Cells.Item(121, 12).Value = 41995
$ws.Cells.Item(121, 14).Value = -45489

# --- BSM sheet (sheet index 3) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(2, 8).Value = 23000
$ws.Cells.Item(2, 10).Value = 23000
$ws.Cells.Item(2, 12).Value = 23000
$ws.Cells.Item(2, 14).Value = -23226
$ws.Cells.Item(6, 8).Value = 23856
$ws.Cells.Item(6, 10).Value = 23856
$ws.Cells.Item(6, 12).Value = 23856
$ws.Cells.Item(6, 14).Value = -24082
$ws.Cells.Item(13, 8).Value = 29000
$ws.Cells.Item(13, 10).Value = 29000
$ws.Cells.Item(13, 12).Value = 29000
$ws.Cells.Item(13, 14).Value = -29336
$ws.Cells.Item(20, 8).Value = 34293.387
$ws.Cells.Item(20, 9).Value = 49539.145
$ws.Cells.Item(20, 10).Value = 2277.3
$ws.Cells.Item(20, 11).Value = 49539.145
$ws.Cells.Item(20, 12).Value = 2277.3
$ws.Cells.Item(20, 13).Value = -49292.145
$ws.Cells.Item(20, 14).Value = -2771.3
$ws.Cells.Item(51, 8).Value = 22776.666
$ws.Cells.Item(51, 10).Value = 22776.666
$ws.Cells.Item(51, 12).Value = 22776.666
$ws.Cells.Item(51, 14).Value = -23758.666
$ws.Cells.Item(53, 8).Value = 24500
$ws.Cells.Item(53, 10).Value = 24500
$ws.Cells.Item(53, 12).Value = 24500
$ws.Cells.Item(53, 14).Value = -25648
$ws.Cells.Item(116, 8).Value = 21496.8
$ws.Cells.Item(116, 10).Value = 21496.8
$ws.Cells.Item(116, 12).Value = 21496.8
$ws.Cells.Item(116, 14).Value = -30674.8
$ws.Cells.Item(117, 8).Value = 50742
$ws.Cells.Item(117, 10).Value = 50742
$ws.Cells.Item(117, 12).Value = 50742
$ws.Cells.Item(117, 14).Value = -59920
$ws.Cells.Item(119, 8).Value = 25000
$ws.Cells.Item(119, 10).Value = 25000
$ws.Cells.Item(119, 12).Value = 25000
$ws.Cells.Item(119, 14).Value = -34676
$ws.Cells.Item(120, 8).Value = 35587
$ws.Cells.Item(120, 10).Value = 35587
$ws.Cells.Item(120, 12).Value = 35587
$ws.Cells.Item(120, 14).Value = -45263

# --- CRP sheet (sheet index 4) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3054.8645
$ws.Cells.Item(31, 9).Value = 2324.9688
$ws.Cells.Item(31, 10).Value = 3919.926
$ws.Cells.Item(31, 11).Value = 2324.9688
$ws.Cells.Item(31, 12).Value = 3919.926
$ws.Cells.Item(31, 13).Value = -2029.9688
$ws.Cells.Item(31, 14).Value = -4509.925999999999
$ws.Cells.Item(34, 8).Value = 3054.8645
$ws.Cells.Item(34, 9).Value = 2324.9688
$ws.Cells.Item(34, 10).Value = 3919.926
$ws.Cells.Item(34, 11).Value = 2324.9688
$ws.Cells.Item(34, 12).Value = 3919.926
$ws.Cells.Item(34, 13).Value = -2122.9688
$ws.Cells.Item(34, 14).Value = -4323.925999999999
$ws.Cells.Item(116, 8).Value = 23000
$ws.Cells.Item(116, 10).Value = 23000
$ws.Cells.Item(116, 12).Value = 23000
$ws.Cells.Item(116, 14).Value = -32178
$ws.Cells.Item(118, 8).Value = 24500
$ws.Cells.Item(118, 10).Value = 24500
$ws.Cells.Item(118, 12).Value = 24500
$ws.Cells.Item(118, 14).Value = -27814
$ws.Cells.Item(119, 8).Value = 36380.5
$ws.Cells.Item(119, 10).Value = 36380.5
$ws.Cells.Item(119, 12).Value = 36380.5
$ws.Cells.Item(119, 14).Value = -46056.5
$ws.Cells.Item(120, 8).Value = 58845.832
$ws.Cells.Item(120, 10).Value = 58845.832
$ws.Cells.Item(120, 12).Value = 58845.832
$ws.Cells.Item(120, 14).Value = -66103.83199999999

# --- CUL sheet (sheet index 5) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 4650.387
$ws.Cells.Item(134, 10).Value = 7519.5
$ws.Cells.Item(134, 12).Value = 22558.5
$ws.Cells.Item(134, 14).Value = -32698.5
$ws.Cells.Item(138, 8).Value = 2024.5714
$ws.Cells.Item(138, 9).Value = 1042.5
$ws.Cells.Item(138, 10).Value = 3334
$ws.Cells.Item(138, 11).Value = 3127.5
$ws.Cells.Item(138, 12).Value = 10002
$ws.Cells.Item(138, 13).Value = 2012.5
$ws.Cells.Item(138, 14).Value = -20282

# --- GSM sheet (sheet index 6) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6332.826
$ws.Cells.Item(70, 9).Value = 5909.4287
$ws.Cells.Item(70, 10).Value = 6518.0625
$ws.Cells.Item(70, 11).Value = 5909.4287
$ws.Cells.Item(70, 12).Value = 6518.0625
$ws.Cells.Item(70, 13).Value = -5639.4287
$ws.Cells.Item(70, 14).Value = -7058.0625
$ws.Cells.Item(73, 8).Value = 6332.826
$ws.Cells.Item(73, 9).Value = 5909.4287
$ws.Cells.Item(73, 10).Value = 6518.0625
$ws.Cells.Item(73, 11).Value = 5909.4287
$ws.Cells.Item(73, 12).Value = 6518.0625
$ws.Cells.Item(73, 13).Value = -4973.4287
$ws.Cells.Item(73, 14).Value = -8390.0625
$ws.Cells.Item(136, 8).Value = 18477.25
$ws.Cells.Item(136, 10).Value = 18477.25
$ws.Cells.Item(136, 12).Value = 55431.75
$ws.Cells.Item(136, 14).Value = -60531.75

# --- LTW sheet (sheet index 7) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 1848.16
$ws.Cells.Item(136, 9).Value = 1729.6136
$ws.Cells.Item(136, 10).Value = 2717.5
$ws.Cells.Item(136, 11).Value = 5188.8408
$ws.Cells.Item(136, 12).Value = 8152.5
$ws.Cells.Item(136, 13).Value = -2638.8408
$ws.Cells.Item(136, 14).Value = -13252.5
$ws.Cells.Item(140, 8).Value = 68857.22
$ws.Cells.Item(140, 10).Value = 68857.22
$ws.Cells.Item(140, 12).Value = 68857.22
$ws.Cells.Item(140, 14).Value = -79217.22

# --- WVR sheet (sheet index 8) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(116, 8).Value = 40914.285
$ws.Cells.Item(116, 10).Value = 40914.285
$ws.Cells.Item(116, 12).Value = 40914.285
$ws.Cells.Item(116, 14).Value = -50092.285
$ws.Cells.Item(117, 8).Value = 31000
$ws.Cells.Item(117, 10).Value = 31000
$ws.Cells.Item(117, 12).Value = 31000
$ws.Cells.Item(117, 14).Value = -40178
$ws.Cells.Item(118, 8).Value = 28473
$ws.Cells.Item(118, 10).Value = 28473
$ws.Cells.Item(118, 12).Value = 28473
